$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 with new computed values
$ws.Range("B2").Value = 0.5127034053968609
$ws.Range("C2").Value = 2.354683216389064
$ws.Range("D2").Value = 21.51883753708786
$ws.Range("E2").Value = 4.63884010686808
$ws.Range("F2").Value = 4.656295975374199
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.3819122929006201
$ws.Range("C3").Value = 2.700847678713336
$ws.Range("D3").Value = 23.37810112586181
$ws.Range("E3").Value = 4.835090601618734
$ws.Range("F3").Value = 4.868918931091797
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.6176629386333392
$ws.Range("C4").Value = 2.723463983124111
$ws.Range("D4").Value = 22.08313540798907
$ws.Range("E4").Value = 4.699269667511013
$ws.Range("F4").Value = 4.706776513729007
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.381006056814591
$ws.Range("C5").Value = 2.954926706810196
$ws.Range("D5").Value = 28.55326029953887
$ws.Range("E5").Value = 5.343525081773161
$ws.Range("F5").Value = 5.386327341582761
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.6469176831876277
$ws.Range("C6").Value = 2.826642545086848
$ws.Range("D6").Value = 25.35880009863441
$ws.Range("E6").Value = 5.035752187968984
$ws.Range("F6").Value = 5.048017233789481
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.2154064431303112
$ws.Range("C7").Value = 2.800487624318238
$ws.Range("D7").Value = 23.5782612740286
$ws.Range("E7").Value = 4.855745182155732
$ws.Range("F7").Value = 4.90456844077524
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = 0.4014956228238964
$ws.Range("C8").Value = 2.47465533314814
$ws.Range("D8").Value = 22.1750436314307
$ws.Range("E8").Value = 4.709038503923142
$ws.Range("F8").Value = 4.744908803445824
$ws.Range("G8").Value = 45

$ws.Range("B9").Value = 0.2593136598854769
$ws.Range("C9").Value = 2.81160301979944
$ws.Range("D9").Value = 23.53166153297062
$ws.Range("E9").Value = 4.850944395988334
$ws.Range("F9").Value = 4.90001041808519
$ws.Range("G9").Value = 44

$ws.Range("B10").Value = 0.5085583818752069
$ws.Range("C10").Value = 2.833416566355677
$ws.Range("D10").Value = 27.39546393872693
$ws.Range("E10").Value = 5.234067628405935
$ws.Range("F10").Value = 5.270953174329484
$ws.Range("G10").Value = 43

# New row 11 for Q9
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.2444862114304603
$ws.Range("C11").Value = 2.496374826965832
$ws.Range("D11").Value = 22.44915894541807
$ws.Range("E11").Value = 4.738054341754436
$ws.Range("F11").Value = 4.789098872386438
$ws.Range("G11").Value = 42
